# The deck's applied design theme is changed from "Integral" back to the
# default "Office Theme": the 12 theme colours used by the slide master
# (ppt/theme/theme1.xml) are updated from the Integral palette to the
# standard Office palette, and the theme / colour-scheme display names are
# updated to match.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

function Set-ThemeRGB($ColorScheme, $Index, $Hex) {
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    $ColorScheme.Item($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeRGB $colorScheme 1  "000000"
Set-ThemeRGB $colorScheme 2  "FFFFFF"
Set-ThemeRGB $colorScheme 3  "44546A"
Set-ThemeRGB $colorScheme 4  "E7E6E6"
Set-ThemeRGB $colorScheme 5  "5B9BD5"
Set-ThemeRGB $colorScheme 6  "ED7D31"
Set-ThemeRGB $colorScheme 7  "A5A5A5"
Set-ThemeRGB $colorScheme 8  "FFC000"
Set-ThemeRGB $colorScheme 9  "4472C4"
Set-ThemeRGB $colorScheme 10 "70AD47"
Set-ThemeRGB $colorScheme 11 "0563C1"
Set-ThemeRGB $colorScheme 12 "954F72"

# Rename the theme / colour scheme to match the stock "Office Theme".
$theme.Name = "Office Theme"
$colorScheme.Name = "Office"
